$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.327.29"
$ws.Cells.Item(2, 5).Value = "  +1.51%  "
$ws.Cells.Item(3, 4).Value = "3.414.99"
$ws.Cells.Item(3, 5).Value = "  +3.47%  "
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "407.29"
$ws.Cells.Item(5, 5).Value = "  +0.57%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "128.47"
$ws.Cells.Item(6, 5).Value = "  +16.49%  "
$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.601"
$ws.Cells.Item(7, 5).Value = "  +7.82%  "
$ws.Cells.Item(8, 2).Value = "LidoStakedEther"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Cells.Item(8, 4).Value = "3.407.29"
$ws.Cells.Item(8, 5).Value = "  +3.47%  "
$ws.Cells.Item(9, 5).Value = "  +0.16%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.673"
$ws.Cells.Item(10, 5).Value = "  +9.79%  "
$ws.Cells.Item(11, 5).Value = "  +18.36%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "42.39"
$ws.Cells.Item(12, 5).Value = "  +10.63%  "
$ws.Cells.Item(13, 5).Value = "  -0.80%  "
$ws.Cells.Item(14, 4).Value = "3.966.79"
$ws.Cells.Item(14, 5).Value = "  +2.85%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "8.59"
$ws.Cells.Item(15, 5).Value = "  +6.34%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "19.73"
$ws.Cells.Item(16, 5).Value = "  +5.20%  "
$ws.Cells.Item(17, 4).Value = "3.430.11"
$ws.Cells.Item(17, 5).Value = "  -1.48%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "11.67"
$ws.Cells.Item(18, 5).Value = "  +13.34%  "
$ws.Cells.Item(19, 4).Value = "61.527.09"
$ws.Cells.Item(19, 5).Value = "  +1.76%  "
$ws.Cells.Item(20, 5).Value = "  +3.54%  "
$ws.Cells.Item(21, 5).Value = "  +23.68%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "3.25"
$ws.Cells.Item(22, 5).Value = "  +1.25%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "82.97"
$ws.Cells.Item(23, 5).Value = "  +13.93%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "13.01"
$ws.Cells.Item(24, 5).Value = "  +7.61%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "307.86"
$ws.Cells.Item(25, 5).Value = "  +6.45%  "
$ws.Cells.Item(26, 5).Value = "  +0.40%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "8.63"
$ws.Cells.Item(27, 5).Value = "  +17.37%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "29.71"
$ws.Cells.Item(28, 5).Value = "  +5.15%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "4.54"
$ws.Cells.Item(29, 5).Value = "  +1.36%  "
$ws.Cells.Item(30, 5).Value = "  +4.00%  "
$ws.Cells.Item(31, 5).Value = "  +6.55%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.116"
$ws.Cells.Item(32, 5).Value = "  +6.94%  "
$ws.Cells.Item(33, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "43.09"
$ws.Cells.Item(33, 5).Value = "  +15.10%  "
$ws.Cells.Item(34, 2).Value = "Cosmos"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "11.69"
$ws.Cells.Item(34, 5).Value = "  +5.80%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "2.57"
$ws.Cells.Item(35, 5).Value = "  +8.08%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Cells.Item(36, 5).Value = "  +0.62%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.0485"
$ws.Cells.Item(37, 5).Value = "  +2.80%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "52.14"
$ws.Cells.Item(38, 5).Value = "  +0.59%  "
$ws.Cells.Item(39, 5).Value = "  -0.14%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "3.41"
$ws.Cells.Item(40, 5).Value = "  +5.17%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "2.98"
$ws.Cells.Item(41, 5).Value = "  +0.57%  "
$ws.Cells.Item(42, 5).Value = "  +5.96%  "
$ws.Cells.Item(43, 5).Value = "  +6.71%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "134.66"
$ws.Cells.Item(44, 5).Value = "  -1.68%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.286"
$ws.Cells.Item(45, 5).Value = "  +5.30%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "16.97"
$ws.Cells.Item(46, 5).Value = "  +6.70%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "3.91"
$ws.Cells.Item(47, 5).Value = "  +5.32%  "
$ws.Cells.Item(48, 5).Value = "  +0.70%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "21.79"
$ws.Cells.Item(49, 5).Value = "  -21.14%  "
$ws.Cells.Item(50, 4).Value = "3.766.69"
$ws.Cells.Item(50, 5).Value = "  -8.21%  "
$ws.Cells.Item(51, 4).Value = "2.152.68"
$ws.Cells.Item(51, 5).Value = "  +1.18%  "
